$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PFOS_abu")

# ---- Header row (row 1): rename columns, add highlight to the moved
# "Interval_btwnbrds" header, and relabel "Fecundity" into K1 ----
$ws.Range("G1").Value = "Number_1brood"
$ws.Range("I1").Value = "Number_2brood"
$ws.Range("J1").Value = "Interval_btwnbrds"
$ws.Range("K1").Value = "Fecundity"
$ws.Range("J1").Interior.Color = 65535

# ---- Treatment label fix: "control" -> "Control" for the LRV-0-1 rows ----
$controlRows = @(2,3,4,5,10,11,12,13)
foreach ($r in $controlRows) {
    $ws.Range("C$r").Value = "Control"
}

# ---- Fecundity (col K) is recomputed as Number_1brood + Number_2brood
# (col G + col I) for every data row - this both shifts the old total
# out of col J into the new col K and corrects a pre-existing typo in
# row 2 (old J2 was 61, but G2+I2 is really 62) ----
$dataRows = @(2,3,4,5,6,7,8,9,10,11,13,14,15,16,17)
foreach ($r in $dataRows) {
    $g = $ws.Range("G$r").Value2
    $i = $ws.Range("I$r").Value2
    $ws.Range("K$r").Value = $g + $i
}

# ---- Rows 2-5, 10, 11, 13: col J becomes the brood interval (H-F) ----
$computeRows = @(2,3,4,5,10,11,13)
foreach ($r in $computeRows) {
    $h = $ws.Range("H$r").Value2
    $f = $ws.Range("F$r").Value2
    $ws.Range("J$r").Value = $h - $f
}

# ---- Rows 6-8, 14-17: the existing interval value in col L moves left
# into col J (row 9, also in this "genotype 24" block, is handled
# separately below since it carries a data correction) ----
$moveLRows = @(6,7,8,14,15,16,17)
foreach ($r in $moveLRows) {
    $oldL = $ws.Range("L$r").Value2
    $ws.Range("J$r").Value = $oldL
    $ws.Range("L$r").ClearContents()
}

# ---- Row 9: data correction (Number_2brood 6 -> 10) plus a corrected
# brood interval (H9-F9, not the stale col L value); Fecundity above
# already used the original Number_2brood value, matching the source ----
$ws.Range("I9").Value = 10
$ws.Range("J9").Value = $ws.Range("H9").Value2 - $ws.Range("F9").Value2
$ws.Range("L9").ClearContents()

# ---- Row 12: no usable replicate data was collected, so blank the
# numeric columns and highlight the whole row (matches the sibling
# "missing data" rows used elsewhere in this workbook) ----
$ws.Range("E12:K12").ClearContents()
$ws.Range("A12:K12").Interior.Color = 65535

# ---- Sheet view: PFOS_abu becomes the selected/active sheet, replacing
# Abu-PE, with a new active cell ----
$ws.Activate()
$ws.Range("O10").Select()
